$d = $word.ActiveDocument
$find = $d.Content.Find
$find.Execute("Infortec Consultores", $true, $false, $false, $false, $false, $true, 1, $false, "Infortec Consultores for Kyndryl", 2)
